# Added Week 15 simulations
$wb = $excel.ActiveWorkbook

# --- OFF sheet updates ---
$wsOff = $wb.Worksheets.Item("OFF")
$wsOff.Range("B2").Value = 256
$wsOff.Range("C2").Value = 174
$wsOff.Range("D2").Value = 54
$wsOff.Range("E2").Value = 20

# --- DEF sheet updates ---
$wsDef = $wb.Worksheets.Item("DEF")
$wsDef.Range("B2").Value = 247
$wsDef.Range("C2").Value = 174
$wsDef.Range("D2").Value = 51
$wsDef.Range("F2").Value = 7
